$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (stored as literal text, matching the
# source workbook convention of text-typed Price/Volume(1h) columns).
$updates = @{
    "D2" = "283.24"
    "E2" = "1.94%"
    "D3" = "28.43"
    "E3" = "4.36%"
    "D4" = "5.054"
    "E4" = "3.72%"
    "D5" = "0.06495"
    "E5" = "0.98%"
    "E6" = "3.24%"
    "D7" = "1.412"
    "E7" = "18.78%"
    "D8" = "0.9178"
    "E8" = "3.97%"
    "D9" = "0.1537"
    "E9" = "-1.43%"
    "D10" = "0.06645"
    "E10" = "29.99%"
    "D11" = "0.07615"
    "E11" = "2.45%"
    "D12" = "0.02783"
    "E12" = "-3.62%"
    "D13" = "0.08970"
    "E13" = "-0.09%"
    "D14" = "0.001589"
    "E14" = "1.22%"
    "D15" = "0.0006336"
    "E15" = "-0.98%"
    "D16" = "0.006150"
    "E16" = "-0.06%"
    "D17" = "3.451"
    "E17" = "-0.92%"
    "E18" = "1.59%"
    "D19" = "2.242"
    "E19" = "-1.43%"
    "D20" = "0.3181"
    "E20" = "-0.08%"
    "E21" = "-0.61%"
    "D22" = "4.007"
    "E22" = "2.20%"
    "E23" = "2.95%"
    "E24" = "0.71%"
    "E25" = "0.43%"
    "D26" = "0.004459"
    "E26" = "14.99%"
    "E27" = "0.00%"
    "E28" = "1.72%"
    "D29" = "0.0001619"
    "E29" = "-1.88%"
    "D40" = "0.04115"
    "E40" = "-0.75%"
    "D41" = "0.006657"
    "E41" = "-2.01%"
    "D42" = "0.1234"
    "E42" = "5.10%"
    "D43" = "0.002050"
    "E43" = "6.79%"
    "D44" = "0.01156"
    "E44" = "-1.77%"
    "D45" = "0.00005395"
    "E45" = "1.70%"
    "E46" = "-0.05%"
    "D47" = "1.933"
    "E47" = "14.73%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text interpretation so numeric-looking strings (and the
    # percentages) round-trip as literal text, not as Number/Percent values.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Drop the temporary text format so the cell keeps its original
    # (unstyled) appearance once the text value has been committed.
    $cell.ClearFormats()
}
